# feat: add 2022-Q1 data
#
# The workbook tracks one worksheet per quarter (fund holdings detail) plus a
# rolling "总计" (Total) overview sheet. Adding a new quarter means:
#   1) the *existing* "总计" sheet (which held the rolled-up overview table)
#      is repurposed into the new quarter's detail sheet "2022-Q1" - its
#      name and content change, but it keeps its original sheetId/r:id slot;
#   2) a brand-new "总计" sheet is appended after it, carrying the same
#      overview table as before plus one new leading row for "2022-Q1".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value as literal TEXT even when it looks like a number
# (mirrors typing NumberFormat "@" then the digits into Excel - the stored
# cell keeps its general/default style afterwards).
# ---------------------------------------------------------------------------
function Set-TextValue($sheet, $row, $col, $text) {
    $cell = $sheet.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------------
# Helper: write a plain number.
# ---------------------------------------------------------------------------
function Set-NumberValue($sheet, $row, $col, $number) {
    $sheet.Cells.Item($row, $col).Value = $number
}

# ---------------------------------------------------------------------------
# Helper: apply the bold/centered/bordered "header" style already used
# elsewhere in this workbook by copying it from a cell that already has it
# (keeps reusing the same cellXf instead of minting a new one).
# ---------------------------------------------------------------------------
function Copy-HeaderStyle($styleSourceCell, $destSheet, $row, $col) {
    $styleSourceCell.Copy()
    $destSheet.Cells.Item($row, $col).PasteSpecial(-4122)
}

# A cell that already carries the workbook's header style (bold, bordered,
# center/top aligned) - used as the formatting donor for new header cells.
$headerStyleDonor = $wb.Worksheets.Item("2021-Q4").Range("B1")
# A cell that already carries the workbook's row-index style (bold, bordered,
# center/top aligned, applied to column A).
$indexStyleDonor = $wb.Worksheets.Item("2021-Q4").Range("A2")

# ===========================================================================
# Step 1: turn the current "总计" sheet into the "2022-Q1" detail sheet.
# ===========================================================================
$fundSheet = $wb.Worksheets.Item("总计")
$fundSheet.Cells.Clear()
$fundSheet.Name = "2022-Q1"

# Header row.
Copy-HeaderStyle $headerStyleDonor $fundSheet 1 2
$fundSheet.Cells.Item(1, 2).Value = "基金代码"
Copy-HeaderStyle $headerStyleDonor $fundSheet 1 3
$fundSheet.Cells.Item(1, 3).Value = "基金名称"
Copy-HeaderStyle $headerStyleDonor $fundSheet 1 4
$fundSheet.Cells.Item(1, 4).Value = "基金规模"
Copy-HeaderStyle $headerStyleDonor $fundSheet 1 5
$fundSheet.Cells.Item(1, 5).Value = "股票总仓位"
Copy-HeaderStyle $headerStyleDonor $fundSheet 1 6
$fundSheet.Cells.Item(1, 6).Value = "仓位占比"
Copy-HeaderStyle $headerStyleDonor $fundSheet 1 7
$fundSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
Copy-HeaderStyle $headerStyleDonor $fundSheet 1 8
$fundSheet.Cells.Item(1, 8).Value = "仓位排名"

function Set-FundRow($sheet, $row, $idx, $code, $name, $scale, $position, $ratio, $marketValue, $rank) {
    Copy-HeaderStyle $indexStyleDonor $sheet $row 1
    Set-NumberValue $sheet $row 1 $idx
    Set-TextValue $sheet $row 2 $code
    $sheet.Cells.Item($row, 3).Value = $name
    Set-TextValue $sheet $row 4 $scale
    Set-TextValue $sheet $row 5 $position
    Set-TextValue $sheet $row 6 $ratio
    Set-TextValue $sheet $row 7 $marketValue
    Set-NumberValue $sheet $row 8 $rank
}

Set-FundRow $fundSheet 2 0 "910004" "东方红启恒三年持有期混合型证券投资基金A" "118.66" "86.48" "6.15" "7.2976" 5
Set-FundRow $fundSheet 3 1 "011724" "东方红启恒三年持有期混合型证券投资基金B" "110.24" "86.48" "6.15" "6.7798" 5
Set-FundRow $fundSheet 4 2 "009014" "泓德睿泽混合" "85.92" "93.52" "3.38" "2.9041" 7
Set-FundRow $fundSheet 5 3 "010059" "东方红鼎元3个月定期开放混合" "27.56" "89.58" "6.68" "1.8410" 6
Set-FundRow $fundSheet 6 4 "169105" "东方红睿华沪港深灵活配置混合（LOF）" "20.95" "86.61" "8.12" "1.7011" 2
Set-FundRow $fundSheet 7 5 "070001" "嘉实成长收益混合A" "24.07" "72.54" "4.80" "1.1554" 6
Set-FundRow $fundSheet 8 6 "001705" "泓德战略转型股票" "22.82" "93.94" "4.25" "0.9698" 6
Set-FundRow $fundSheet 9 7 "910009" "东方红启程三年持有期混合型证券投资基金A" "8.06" "87.95" "6.16" "0.4965" 5
Set-FundRow $fundSheet 10 8 "008150" "嘉实远见企业精选两年持有期混合" "14.46" "91.85" "3.37" "0.4873" 8
Set-FundRow $fundSheet 11 9 "002989" "融通通乾研究精选灵活配置混合" "6.13" "94.84" "6.00" "0.3678" 5
Set-FundRow $fundSheet 12 10 "002846" "泓德泓华灵活配置混合" "7.01" "93.24" "3.87" "0.2713" 7
Set-FundRow $fundSheet 13 11 "011011" "融通产业趋势精选2年封闭运作混合" "3.07" "94.85" "5.93" "0.1821" 5
Set-FundRow $fundSheet 14 12 "000870" "嘉实新收益灵活配置混合" "4.17" "91.30" "4.15" "0.1731" 8
Set-FundRow $fundSheet 15 13 "008382" "融通产业趋势股票" "2.39" "94.68" "5.96" "0.1424" 7
Set-FundRow $fundSheet 16 14 "002801" "泓德泓信灵活配置混合" "2.53" "92.18" "3.52" "0.0891" 7
Set-FundRow $fundSheet 17 15 "008110" "九泰科盈价值灵活配置混合A" "3.52" "34.14" "1.53" "0.0539" 8
Set-FundRow $fundSheet 18 16 "011727" "工银瑞信聚瑞混合型证券投资基金A" "3.46" "29.56" "1.28" "0.0443" 8
Set-FundRow $fundSheet 19 17 "009015" "泓德睿享一年持有期混合A" "3.41" "24.56" "1.25" "0.0426" 3
Set-FundRow $fundSheet 20 18 "008136" "九泰科盈价值灵活配置混合C" "2.45" "34.14" "1.53" "0.0375" 8
Set-FundRow $fundSheet 21 19 "910011" "东方红启瑞三年持有混合A" "0.49" "48.85" "1.98" "0.0097" 10
Set-FundRow $fundSheet 22 20 "006603" "嘉实互融精选股票" "0.14" "85.91" "4.14" "0.0058" 7
Set-FundRow $fundSheet 23 21 "011728" "工银瑞信聚瑞混合型证券投资基金C" "0.17" "29.56" "1.28" "0.0022" 8
Set-FundRow $fundSheet 24 22 "011312" "东方红启瑞三年持有混合B" "0.11" "48.85" "1.98" "0.0022" 10
Set-FundRow $fundSheet 25 23 "009016" "泓德睿享一年持有期混合C" "0.07" "24.56" "1.25" "0.0009" 3
Set-FundRow $fundSheet 26 24 "960024" "嘉实成长收益混合H" "0.01" "72.54" "4.80" "0.0005" 6

# ===========================================================================
# Step 2: append a brand-new "总计" overview sheet after "2022-Q1".
# ===========================================================================
$totalSheet = $wb.Worksheets.Add($null, $fundSheet)
$totalSheet.Name = "总计"

Copy-HeaderStyle $headerStyleDonor $totalSheet 1 2
$totalSheet.Cells.Item(1, 2).Value = "日期"
Copy-HeaderStyle $headerStyleDonor $totalSheet 1 3
$totalSheet.Cells.Item(1, 3).Value = "持有数量(只)"
Copy-HeaderStyle $headerStyleDonor $totalSheet 1 4
$totalSheet.Cells.Item(1, 4).Value = "持有市值(亿元)"

function Set-TotalRow($sheet, $row, $idx, $quarter, $count, $marketValue) {
    Copy-HeaderStyle $indexStyleDonor $sheet $row 1
    Set-NumberValue $sheet $row 1 $idx
    $sheet.Cells.Item($row, 2).Value = $quarter
    Set-NumberValue $sheet $row 3 $count
    Set-NumberValue $sheet $row 4 $marketValue
}

Set-TotalRow $totalSheet 2 0 "2022-Q1" 25 25.06
Set-TotalRow $totalSheet 3 1 "2021-Q4" 43 21.8
Set-TotalRow $totalSheet 4 2 "2021-Q3" 57 27.04
Set-TotalRow $totalSheet 5 3 "2021-Q2" 52 21.91
Set-TotalRow $totalSheet 6 4 "2021-Q1" 87 24
Set-TotalRow $totalSheet 7 5 "2020-Q4" 159 74.69
